$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some Price (column D) values are purely numeric-looking strings (e.g. "519.42").
# Excel's COM automation auto-converts such text into a real number when assigned
# via .Value, which would lose the original text formatting. Force those specific
# cells to Text format first so the assigned strings are preserved verbatim.
$textForcedCells = @("D5","D6","D7","D12","D15","D19","D20","D23","D25","D31","D33","D39","D40","D42","D43","D44","D46","D48")
foreach ($addr in $textForcedCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "58.123.95"
$ws.Range("E2").Value = "  -1.29%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.471.85"
$ws.Range("E3").Value = "  -1.97%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.15%  "

# Row 5 - BNB
$ws.Range("D5").Value = "519.42"
$ws.Range("E5").Value = "  -3.12%  "

# Row 6 - Solana
$ws.Range("D6").Value = "132.30"
$ws.Range("E6").Value = "  -3.95%  "

# Row 7 - USDC
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.04%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  -1.70%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  -2.07%  "

# Row 10 - TRON
$ws.Range("E10").Value = "  -0.81%  "

# Row 11 - Toncoin
$ws.Range("E11").Value = "  +0.07%  "

# Row 12 - Cardano
$ws.Range("D12").Value = "0.342"
$ws.Range("E12").Value = "  -1.77%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "2.910.26"
$ws.Range("E13").Value = "  -1.96%  "

# Row 14 - WrappedBTC
$ws.Range("D14").Value = "58.061.58"
$ws.Range("E14").Value = "  -1.41%  "

# Row 15 - Avalanche
$ws.Range("D15").Value = "22.03"
$ws.Range("E15").Value = "  -4.49%  "

# Row 16 - ShibaInu
$ws.Range("E16").Value = "  -2.14%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "2.473.65"
$ws.Range("E17").Value = "  -1.83%  "

# Row 18 - Chainlink
$ws.Range("E18").Value = "  -2.47%  "

# Row 19 - becomes Polkadot (was BitcoinCash)
$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D19").Value = "4.18"
$ws.Range("E19").Value = "  -2.68%  "

# Row 20 - becomes BitcoinCash (was Polkadot)
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "319.40"
$ws.Range("E20").Value = "  -1.93%  "

# Row 21 - Dai
$ws.Range("E21").Value = "  -0.06%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  -3.40%  "

# Row 23 - Litecoin
$ws.Range("D23").Value = "64.26"
$ws.Range("E23").Value = "  -2.39%  "

# Row 24 - Polygon
$ws.Range("E24").Value = "  -3.70%  "

# Row 25 - Binance-PegBSC-USD
$ws.Range("D25").Value = "0.999"
$ws.Range("E25").Value = "  -0.23%  "

# Row 26 - Kaspa
$ws.Range("E26").Value = "  -3.71%  "

# Row 27 - InternetComputer(DFINITY)
$ws.Range("E27").Value = "  -3.41%  "

# Row 28 - PEPE
$ws.Range("D28").Value = "0.0₃0750"
$ws.Range("E28").Value = "  -2.75%  "

# Row 29 - Aptos
$ws.Range("E29").Value = "  -5.22%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  -4.90%  "

# Row 31 - Monero
$ws.Range("D31").Value = "166.65"
$ws.Range("E31").Value = "  +2.71%  "

# Row 32 - Fetch.AI
$ws.Range("E32").Value = "  -3.92%  "

# Row 33 - USDe
$ws.Range("D33").Value = "0.999"
$ws.Range("E33").Value = "  +0.00%  "

# Row 34 - FirstDigitalUSD
$ws.Range("E34").Value = "  +0.09%  "

# Row 35 - EthereumClassic
$ws.Range("E35").Value = "  -1.89%  "

# Row 36 - ImmutableX
$ws.Range("E36").Value = "  -10.38%  "

# Row 37 - NEARProtocol
$ws.Range("E37").Value = "  -3.27%  "

# Row 38 - Stacks
$ws.Range("E38").Value = "  -4.00%  "

# Row 39 - SuiNetwork
$ws.Range("D39").Value = "0.795"
$ws.Range("E39").Value = "  -2.87%  "

# Row 40 - Bittensor
$ws.Range("D40").Value = "276.91"
$ws.Range("E40").Value = "  -3.35%  "

# Row 41 - Filecoin
$ws.Range("E41").Value = "  -4.83%  "

# Row 42 - RenderToken
$ws.Range("D42").Value = "5.04"
$ws.Range("E42").Value = "  -2.87%  "

# Row 43 - Mantle
$ws.Range("D43").Value = "0.596"
$ws.Range("E43").Value = "  -2.32%  "

# Row 44 - Aave
$ws.Range("D44").Value = "125.99"
$ws.Range("E44").Value = "  -4.86%  "

# Row 45 - Stellar
$ws.Range("E45").Value = "  -2.55%  "

# Row 46 - Hedera
$ws.Range("D46").Value = "0.0492"
$ws.Range("E46").Value = "  -3.56%  "

# Row 47 - VeChain
$ws.Range("E47").Value = "  -3.51%  "

# Row 48 - InjectiveProtocol
$ws.Range("D48").Value = "17.11"
$ws.Range("E48").Value = "  -1.54%  "

# Row 49 - Maker
$ws.Range("D49").Value = "1.736.24"
$ws.Range("E49").Value = "  -1.45%  "

# Row 50 - BitgetToken
$ws.Range("E50").Value = "  -1.78%  "

# Row 51 - ZEEBU
$ws.Range("E51").Value = "  -1.76%  "
